$wb = $excel.ActiveWorkbook

# --- Settings sheet ---
$wsSettings = $wb.Worksheets.Item("Settings")

# Row 2: rename queue asset from "ProcessABCQueue" to "Queue"
$wsSettings.Range("B2").Value = "Queue"

# --- Assets sheet ---
$wsAssets = $wb.Worksheets.Item("Assets")

# New rows for tempPath / outputPath assets
$wsAssets.Range("A2").Value = "tempPath"
$wsAssets.Range("B2").Value = "tempPath"
$wsAssets.Range("A3").Value = "outputPath"
$wsAssets.Range("B3").Value = "outputPath"

# Move the selection/active-cell on Assets off of the header selection,
# and make sure Assets is no longer the tab-selected sheet (Settings will be).
$wsAssets.Range("B7").Select()

# --- Finally make Settings the active/selected sheet with its new selection ---
$wsSettings.Activate()
$wsSettings.Range("A10").Select()
